$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "Through 2022-05-12" to "Through 2022-05-13"
$ws.Name = "Through 2022-05-13"

# Update the header cell text (column I, row 1) for the 2022 partial-year total
$ws.Range("I1").Value = "2022 (through 05-13)"

# Update May 2022 total (row 6 = May)
$ws.Range("I6").Value = 48

# Update overall Total row for 2022 (row 14 = Total)
$ws.Range("I14").Value = 600
